# Apply the target edit: insert a new data row at row 9 (pushing existing
# rows 9-55 down to 10-56), then populate the new row 9 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 9; this shifts rows 9..55 -> 10..56
# and Excel automatically extends the used range / dimension to T56.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data.
$ws.Range("A9").Value2 = 9
$ws.Range("B9").Value2 = 'Vega Central Mapocho de Santiago'
$ws.Range("C9").Value2 = 'Metropolitana'
$ws.Range("D9").Value2 = 45050
$ws.Range("E9").Value2 = 13
$ws.Range("F9").Value2 = 'Fruta'
$ws.Range("G9").Value2 = 100107
$ws.Range("H9").Value2 = 'Otros'
$ws.Range("I9").Value2 = 100107001
$ws.Range("J9").Value2 = 'Caqui'
$ws.Range("K9").Value2 = 'Mankaki'
$ws.Range("L9").Value2 = 'Primera'
$ws.Range("M9").Value2 = 400
$ws.Range("N9").Value2 = 9000
$ws.Range("O9").Value2 = 9500
$ws.Range("P9").Value2 = 9225
$ws.Range("Q9").Value2 = '$/caja 12 kilos empedrada'
$ws.Range("R9").Value2 = 'Región del Maule'
$ws.Range("S9").Value2 = 769
$ws.Range("T9").Value2 = 12
